# Split the "Solutions architecture  and big data" run on slide 3
# (Roadmap / Content Placeholder 2) into four runs and remove the
# duplicated space between "architecture" and "and":
#
#   "Solutions architecture  and big data"
# ->  "Solutions " + "architecture " + "and " + "big data"

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(3)
$shp = $s.Shapes.Item(2)
$tr  = $shp.TextFrame.TextRange

$full = $tr.Text
$needle = "Solutions architecture  and big data"
$idx0 = $full.IndexOf($needle)
$base = $idx0 + 1   # TextRange.Characters() is 1-indexed

# Segment lengths as laid out in the ORIGINAL text (includes the extra space):
$len1 = 10  # "Solutions "
$len2 = 14  # "architecture  "  (two spaces)
$len3 = 4   # "and "
$len4 = 8   # "big data"

# Edit right-to-left so earlier segments' absolute offsets stay valid
# while later segments change length.
$seg4 = $tr.Characters($base + $len1 + $len2 + $len3, $len4)
$seg4.Text = "big data"

$seg3 = $tr.Characters($base + $len1 + $len2, $len3)
$seg3.Text = "and "

$seg2 = $tr.Characters($base + $len1, $len2)
$seg2.Text = "architecture "

$seg1 = $tr.Characters($base, $len1)
$seg1.Text = "Solutions "

Write-Host "Result:" $tr.Text
